# Append two new submission rows to the data-collection sheet
# (the first worksheet, "八位序列号收集收集结果yd5").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 109
$ws.Cells.Item(109, 1).Value = [char]0x3000
$ws.Cells.Item(109, 2).Value = 45977.5766319445
$ws.Cells.Item(109, 2).NumberFormat = "yyyy/m/d h:mm:ss;@"
$ws.Cells.Item(109, 3).Value = "4c2cc8c6"
$ws.Cells.Item(109, 4).Value = "1503403546"

# Row 110
$ws.Cells.Item(110, 1).Value = "cyc"
$ws.Cells.Item(110, 2).Value = 45977.8849652778
$ws.Cells.Item(110, 2).NumberFormat = "yyyy/m/d h:mm:ss;@"
$ws.Cells.Item(110, 3).Value = "cdec554d"
$ws.Cells.Item(110, 4).Value = "2648616816"

Write-Host ("Done. New dim: " + $ws.UsedRange.Address())
